$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the cached "datetimeFigureOut" footer date text from 5/4/21 to
#    5/9/21 everywhere it appears: the slide master, every slide layout, and
#    the notes master.
# ---------------------------------------------------------------------------
function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
                if ($shp.TextFrame.TextRange.Text -eq "5/4/21") {
                    $shp.TextFrame.TextRange.Text = "5/9/21"
                }
            }
        }
    }
}

# Slide master footer date.
Update-DateShapes $p.SlideMaster.Shapes

# Every slide layout footer date.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DateShapes $layouts.Item($li).Shapes
}

# Notes master footer date.
Update-DateShapes $p.NotesMaster.Shapes

# ---------------------------------------------------------------------------
# 2) Zero-index the "center_offset" labels in the diagram on slide 1:
#      center_offset(1) -> center_offset(0)
#      center_offset(2) -> center_offset(1)
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $tr = $shp.TextFrame.TextRange
        $full = $tr.Text
        if ($full -eq "center_offset(1)") {
            $sub = $tr.Characters(14, 3)
            $sub.Text = "(0)"
        }
        elseif ($full -eq "center_offset(2)") {
            $sub = $tr.Characters(14, 3)
            $sub.Text = "(1)"
        }
    }
}
